$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the column headers in row 1 from *_old / *_new to
#        *_FV2404 / *_FV2410 (keeping the "diff" column untouched). ---
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"

$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- 2) Turn the data range into an Excel Table ("Table1") spanning
#        A1:U56, picking up the (just renamed) header row. ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U56"), $null, 1)
$lo.Name = "Table1"

# --- 3) Freeze the header row (split below row 1). ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
